# Update Brazil Summary country indicator figures.
#
# These cells hold numeric-looking values stored as plain text (shared
# strings) in the original workbook. A bare numeric assignment would be
# auto-detected by Excel and coerced into a real number (and pick up
# binary floating-point noise, e.g. 22.95 -> 22.94999999999...), so each
# new value is entered with a leading apostrophe to force it to remain
# text, exactly like typing '22.95 into the cell in the Excel UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

$ws.Range("B11").Value = "'22.95"   # Enterprises density (per 1000 people) - Micro: 22.9 -> 22.95
$ws.Range("D11").Value = "'25.94"   # Enterprises density (per 1000 people) - MSMEs: 25.9 -> 25.94
$ws.Range("D12").Value = "'53.01"   # Employment (% of total) - MSMEs: 53 -> 53.01
$ws.Range("C14").Value = "'11.51"   # Enterprises (% of total) - SMEs: 11.5 -> 11.51
$ws.Range("D14").Value = "'99.61"   # Enterprises (% of total) - MSMEs: 99.6 -> 99.61
